$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: new blog post title + link
$ws.Range("D12").Value = "“구글 브레인 팀에게 배우는 딥러닝 with TensorFlow.js”가 출간되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/03/25/%ea%b5%ac%ea%b8%80-%eb%b8%8c%eb%a0%88%ec%9d%b8-%ed%8c%80%ec%97%90%ea%b2%8c-%eb%b0%b0%ec%9a%b0%eb%8a%94-%eb%94%a5%eb%9f%ac%eb%8b%9d-with-tensorflow-js%ea%b0%80-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88/"

# Row 26: new post title (link unchanged)
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

# Row 36: new seminar title + link
$ws.Range("D36").Value = "Dive into audio transformer"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/359"

# Row 46: new job posting title + link
$ws.Range("D46").Value = "[유한양행] 2022년 03월, 생물정보학(Bioinformatics 채용), 합성신약(AI신약개발) 연구원"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/446"
